$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.8680396643783372
$ws.Range("F2").Value = 0.0157324009208808
$ws.Range("G2").Value = "[0.8700980392156863, 0.868935097668557, 0.8681159420289856, 0.8529411764705882, 0.869296833064949, 0.8509930220075148, 0.855072463768116, 0.8385049580472921, 0.8679633867276888, 0.8840579710144928]"

$ws.Range("E3").Value = 0.855072463768116
$ws.Range("F3").Value = 0.01561167945439046
$ws.Range("G3").Value = "[0.8676470588235294, 0.855072463768116, 0.855072463768116, 0.8529411764705882, 0.8695652173913043, 0.8529411764705882, 0.855072463768116, 0.8405797101449275, 0.8695652173913043, 0.8840579710144928]"

$ws.Range("E4").Value = 0.9048442534908701
$ws.Range("F4").Value = 0.005545112781954797
$ws.Range("G4").Value = "[0.9273684210526315, 0.9, 0.9031578947368422, 0.9194414607948442, 0.9084210526315789, 0.8852631578947369, 0.9042105263157895, 0.9030612244897959, 0.9054779806659505, 0.9087003222341568]"

$ws.Range("E5").Value = 0.8587654828982304
$ws.Range("F5").Value = 0.01379339028942983
$ws.Range("G5").Value = "[0.8686616597471659, 0.8589544513457559, 0.8585765144507049, 0.8529411764705882, 0.8634143399867334, 0.8493343825824645, 0.855072463768116, 0.8393410132540567, 0.8685517381169555, 0.8840579710144928]"

$ws.Range("E6").Value = 0.817037467015322
$ws.Range("F6").Value = 0.01282413995346432
$ws.Range("G6").Value = "[0.8120935330654, 0.7971014492753623, 0.8107585139318886, 0.8235294117647058, 0.8264617691154422, 0.8146684233640756, 0.8681945505474918, 0.8194065106665682, 0.804922935357718, 0.824045407636739]"

$ws.Range("E7").Value = 0.8115942028985508
$ws.Range("F7").Value = 0.01102941176470584
$ws.Range("G7").Value = "[0.8115942028985508, 0.7971014492753623, 0.8088235294117647, 0.8235294117647058, 0.8115942028985508, 0.8115942028985508, 0.855072463768116, 0.8115942028985508, 0.7971014492753623, 0.8235294117647058]"

$ws.Range("E8").Value = 0.8910195707070707
$ws.Range("F8").Value = 0.01587902837902844
$ws.Range("G8").Value = "[0.890625, 0.8897306397306397, 0.9074074074074073, 0.9056277056277056, 0.8914141414141413, 0.9030303030303031, 0.8813131313131313, 0.9097222222222222, 0.8888888888888888, 0.8800505050505051]"

$ws.Range("E9").Value = 0.8115147905499306
$ws.Range("F9").Value = 0.01121040216081437
$ws.Range("G9").Value = "[0.8116734149217173, 0.7971014492753623, 0.8081999923256975, 0.8235294117647058, 0.8114353782013104, 0.8115942028985508, 0.8529194655643316, 0.8114353782013104, 0.7947603121516165, 0.8230678496333145]"

$ws.Range("E10").Value = 0.8776353278421445
$ws.Range("F10").Value = 0.02355354768036888
$ws.Range("G10").Value = "[0.862486126526082, 0.876911530162304, 0.8681626928471249, 0.8783591255219849, 0.8529411764705882, 0.8962414397197006, 0.8783591255219849, 0.8622052610058607, 0.8904920458749538, 0.891475252939568]"

$ws.Range("E11").Value = 0.8823529411764706
$ws.Range("F11").Value = 0.01449275362318847
$ws.Range("G11").Value = "[0.8529411764705882, 0.8823529411764706, 0.8695652173913043, 0.8840579710144928, 0.8529411764705882, 0.8970588235294118, 0.8840579710144928, 0.8695652173913043, 0.8823529411764706, 0.8840579710144928]"

$ws.Range("E12").Value = 0.8961038961038961
$ws.Range("F12").Value = 0.03259240759240767
$ws.Range("G12").Value = "[0.8766233766233766, 0.8853146853146855, 0.9244755244755245, 0.9324675324675326, 0.8961038961038961, 0.9041958041958043, 0.8961038961038961, 0.8831168831168831, 0.8545454545454546, 0.9202797202797204]"

$ws.Range("E13").Value = 0.8764302059496567
$ws.Range("F13").Value = 0.02523296550889165
$ws.Range("G13").Value = "[0.8568316215375039, 0.8785014005602241, 0.8517279821627647, 0.8764302059496567, 0.8529411764705882, 0.8897290485192186, 0.8764302059496567, 0.8634705656021547, 0.8854652972300031, 0.8868491680085883]"
